# "Two more empty superheroes"
# Adds two more (mostly empty) hero rows to the Tabelle1 table/sheet,
# growing the table from A1:E22 to A1:E25, and updates the sheet
# selection to reflect where the editor ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Row 23: a full "empty" hero row (blank name/plus/minus, "but..." divider,
# next svg background image).
$lo.ListRows.Add() | Out-Null
$ws.Cells.Item(23, 1).Value = "  "
$ws.Cells.Item(23, 2).Value = "  "
$ws.Cells.Item(23, 3).Value = "but…"
$ws.Cells.Item(23, 4).Value = "  "
$ws.Cells.Item(23, 5).Value = "svg/l.svg"

# Row 24: only the "but..." divider + background image are set.
$lo.ListRows.Add() | Out-Null
$ws.Cells.Item(24, 3).Value = "but…"
$ws.Cells.Item(24, 5).Value = "svg/w.svg"

# Row 25: same pattern as row 24.
$lo.ListRows.Add() | Out-Null
$ws.Cells.Item(25, 3).Value = "but…"
$ws.Cells.Item(25, 5).Value = "svg/x.svg"

# Reflect the editor's final viewport/selection: scrolled so column C is
# at the left edge, with the cell just past the new last row selected.
$excel.Goto($ws.Range("C1"), $true)
$ws.Range("E26").Select()
